$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 23: fill in previously-empty K/L/M cells (style already set) ---
$ws.Range("K23").Value = 0.39746917586000002
$ws.Range("L23").Value = 0.39740470397400002
$ws.Range("M23").Value = 0.39743693730200003

# --- Row 28: update existing K/L/M values ---
$ws.Range("K28").Value = 0.53282417939600002
$ws.Range("L28").Value = 0.53187347931899998
$ws.Range("M28").Value = 0.53234840490299995

# --- Row 30: add L30 (mirrors J30/H30/F30 formatting + "weighted + avg + cosine" label) ---
$ws.Range("J30").Copy()
$ws.Range("L30").PasteSpecial(-4122)
$ws.Range("L30").Value = "weighted + avg + cosine"

# --- Row 31: add K31/L31/M31, mirroring H31/I31/J31 formatting ---
$ws.Range("H31:J31").Copy()
$ws.Range("K31:M31").PasteSpecial(-4122)
$ws.Range("K31").Value = "Evaluation (TWSI full set)"

# --- Row 32: add K32/L32/M32, mirroring H32/I32/J32 formatting ---
$ws.Range("H32:J32").Copy()
$ws.Range("K32:M32").PasteSpecial(-4122)
$ws.Range("K32").Value = "Precision"
$ws.Range("L32").Value = "Recall"
$ws.Range("M32").Value = "F1"

# --- Row 33: add new data cells K33/L33/M33 (default styling) ---
$ws.Range("K33").Value = 0.66928998201400003
$ws.Range("L33").Value = 0.66914703045299995
$ws.Range("M33").Value = 0.66921849859899996

# --- Row 34: add new data cells K34/L34/M34 (default styling) ---
$ws.Range("K34").Value = 0.62832077955200005
$ws.Range("L34").Value = 0.62818657847600001
$ws.Range("M34").Value = 0.62825367184699998

# --- Update the active selection to match the saved view state ---
$ws.Range("K43").Select()

$excel.CutCopyMode = 0
